$d = $word.ActiveDocument

# The requirement bullet reads:
#   "The application must have at least 5 independent entity models"
# The edit drops the bolded word "independent" (and the space that
# followed it), so the bullet becomes:
#   "The application must have at least 5 entity models"
# and the bold emphasis moves onto "entity models" (it keeps the
# <w:b/> that used to sit on "independent", without picking up the
# <w:bCs/> companion property).

# 1. Remove "independent " (the word plus its trailing space) outright.
#    Doing this with an empty replacement - rather than folding it into
#    a text replace of the whole phrase - lets the now-adjacent,
#    identically-formatted space run merge cleanly instead of carrying
#    "independent"'s run formatting forward.
$null = $d.Content.Find.Execute("independent ", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "", 2)

# 2. Bold the now-isolated "entity models" run on its own, so only
#    <w:b/> is added (no <w:bCs/>, no stray rsid attributes).
$r = $d.Content
$null = $r.Find.Execute("entity models", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$r.Bold = 1
